# Auto-generated edit script applying the Gilgamesh_Profits.xlsx diff
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 90910680
$ws.Range("I18").Value = 1243.2222
$ws.Range("K18").Value = 1243.2222
$ws.Range("M18").Value = -959.2221999999999
$ws.Range("H33").Value = 290.44446
$ws.Range("I33").Value = 314.375
$ws.Range("J33").Value = 99
$ws.Range("K33").Value = 314.375
$ws.Range("L33").Value = 99
$ws.Range("M33").Value = -85.375
$ws.Range("N33").Value = -557
$ws.Range("H88").Value = 16667648
$ws.Range("I88").Value = 33333896
$ws.Range("K88").Value = 33333896
$ws.Range("M88").Value = -33333490
$ws.Range("H91").Value = 16667648
$ws.Range("I91").Value = 33333896
$ws.Range("K91").Value = 33333896
$ws.Range("M91").Value = -33332492
$ws.Range("H94").Value = 39998
$ws.Range("I94").Value = 39998
$ws.Range("K94").Value = 39998
$ws.Range("M94").Value = -39547
$ws.Range("H121").Value = 2248.25
$ws.Range("J121").Value = 2248.25
$ws.Range("L121").Value = 6744.75
$ws.Range("N121").Value = -10238.75
$ws.Range("H138").Value = 327691.1
$ws.Range("I138").Value = 3417.111
$ws.Range("K138").Value = 10251.333
$ws.Range("M138").Value = -5111.332999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1171.7693
$ws.Range("I2").Value = 1062.1428
$ws.Range("J2").Value = 1299.6666
$ws.Range("K2").Value = 1062.1428
$ws.Range("L2").Value = 1299.6666
$ws.Range("M2").Value = -949.1428000000001
$ws.Range("N2").Value = -1525.6666
$ws.Range("H45").Value = 39460.418
$ws.Range("I45").Value = 45752.2
$ws.Range("K45").Value = 45752.2
$ws.Range("M45").Value = -45375.2
$ws.Range("H102").Value = 5586.923
$ws.Range("I102").Value = 5430.091
$ws.Range("K102").Value = 5430.091
$ws.Range("M102").Value = -3808.091
$ws.Range("H110").Value = 2488.8965
$ws.Range("I110").Value = 1450.3182
$ws.Range("K110").Value = 1450.3182
$ws.Range("M110").Value = 594.6818000000001
$ws.Range("H116").Value = 1171.7693
$ws.Range("I116").Value = 1062.1428
$ws.Range("J116").Value = 1299.6666
$ws.Range("K116").Value = 1062.1428
$ws.Range("L116").Value = 1299.6666
$ws.Range("M116").Value = 1231.8572
$ws.Range("N116").Value = -5887.6666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1171.7693
$ws.Range("I3").Value = 1062.1428
$ws.Range("J3").Value = 1299.6666
$ws.Range("K3").Value = 1062.1428
$ws.Range("L3").Value = 1299.6666
$ws.Range("M3").Value = -948.1428000000001
$ws.Range("N3").Value = -1527.6666
$ws.Range("H82").Value = 62730.5
$ws.Range("I82").Value = 35249.25
$ws.Range("J82").Value = 90211.75
$ws.Range("K82").Value = 35249.25
$ws.Range("L82").Value = 90211.75
$ws.Range("M82").Value = -34866.25
$ws.Range("N82").Value = -90977.75
$ws.Range("H85").Value = 62730.5
$ws.Range("I85").Value = 35249.25
$ws.Range("J85").Value = 90211.75
$ws.Range("K85").Value = 35249.25
$ws.Range("L85").Value = 90211.75
$ws.Range("M85").Value = -33923.25
$ws.Range("N85").Value = -92863.75
$ws.Range("H94").Value = 181819310
$ws.Range("I94").Value = 200000240
$ws.Range("K94").Value = 200000240
$ws.Range("M94").Value = -199999789
$ws.Range("H95").Value = 80778
$ws.Range("J95").Value = 80778
$ws.Range("L95").Value = 80778
$ws.Range("N95").Value = -86270

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.888885
$ws.Range("I7").Value = 84.42856999999999
$ws.Range("J7").Value = 131.5
$ws.Range("K7").Value = 84.42856999999999
$ws.Range("L7").Value = 131.5
$ws.Range("M7").Value = 28.57143000000001
$ws.Range("N7").Value = -357.5
$ws.Range("H92").Value = 32297.2
$ws.Range("J92").Value = 32297.2
$ws.Range("L92").Value = 32297.2
$ws.Range("N92").Value = -37289.2
$ws.Range("H107").Value = 689.4
$ws.Range("I107").Value = 689.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 689.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1230.6
$ws.Range("N107").Value = $null
$ws.Range("H122").Value = 3469.6155
$ws.Range("I122").Value = 1858.875
$ws.Range("J122").Value = 6046.8
$ws.Range("K122").Value = 5576.625
$ws.Range("L122").Value = 18140.4
$ws.Range("M122").Value = -3126.625
$ws.Range("N122").Value = -23040.4
$ws.Range("H132").Value = 2330.3076
$ws.Range("I132").Value = 2087.3333
$ws.Range("J132").Value = 3666.6667
$ws.Range("K132").Value = 6261.999899999999
$ws.Range("L132").Value = 11000.0001
$ws.Range("M132").Value = -3731.999899999999
$ws.Range("N132").Value = -16060.0001
$ws.Range("H134").Value = 2964.1785
$ws.Range("I134").Value = 2630.5217
$ws.Range("J134").Value = 4499
$ws.Range("K134").Value = 7891.5651
$ws.Range("L134").Value = 13497
$ws.Range("M134").Value = -5356.5651
$ws.Range("N134").Value = -18567

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 179
$ws.Range("I28").Value = 179
$ws.Range("K28").Value = 537
$ws.Range("M28").Value = -305
$ws.Range("H48").Value = 4500
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = $null
$ws.Range("H141").Value = 11680.25
$ws.Range("I141").Value = 5144.923
$ws.Range("K141").Value = 15434.769
$ws.Range("M141").Value = -10254.769

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 2044.8
$ws.Range("J36").Value = 2044.8
$ws.Range("L36").Value = 2044.8
$ws.Range("N36").Value = -3014.8
$ws.Range("H40").Value = 35000
$ws.Range("J40").Value = 35000
$ws.Range("L40").Value = 35000
$ws.Range("N40").Value = -35302
$ws.Range("H112").Value = 90000
$ws.Range("I112").Value = 60000
$ws.Range("J112").Value = 100000
$ws.Range("K112").Value = 60000
$ws.Range("L112").Value = 100000
$ws.Range("M112").Value = -58892
$ws.Range("N112").Value = -102216
$ws.Range("H139").Value = 81359.10000000001
$ws.Range("J139").Value = 82588.336
$ws.Range("L139").Value = 82588.336
$ws.Range("N139").Value = -92868.336

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 29435.572
$ws.Range("I56").Value = 22409.8
$ws.Range("J56").Value = 47000
$ws.Range("K56").Value = 22409.8
$ws.Range("L56").Value = 47000
$ws.Range("M56").Value = -21718.8
$ws.Range("N56").Value = -48382

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 382.85184
$ws.Range("J113").Value = 326.55554
$ws.Range("L113").Value = 979.66662
$ws.Range("N113").Value = -5319.66662
